{"js": "// Remove the 4 paragraphs that immediately follow the\n// \"LOQ4031: Qu\u00edmica Geral I (Requisito fraco)\" paragraph:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. an empty paragraph\n//   4. an empty paragraph with a page-break-before\n// The two paragraphs that come after them (another empty paragraph and a\n// final page-break-before empty paragraph) are left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its text.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"LOQ4031: Qu\u00edmica Geral I (Requisito fraco)\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4031' paragraph\");\n}\n\n// Delete the 4 paragraphs that follow the anchor. Deleting from the\n// last one back to the first keeps the remaining indices valid.\nconst toDeleteCount = 4;\nfor (let i = anchorIndex + toDeleteCount; i > anchorIndex; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the 4 paragraphs that immediately follow the\n# \"LOQ4031: Qu\u00edmica Geral I (Requisito fraco)\" paragraph:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. an empty paragraph\n#   4. an empty paragraph with a page-break-before\n# The two paragraphs that come after them (another empty paragraph and a\n# final page-break-before empty paragraph) are left untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph via Find so we do not depend on a hard-coded\n# paragraph index.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOQ4031: Qu\u00edmica Geral I (Requisito fraco)\")\nif (-not $found) {\n    throw \"Could not find the 'LOQ4031' anchor paragraph\"\n}\n[void]$rng.Expand(4)  # wdParagraph -> grow the match to the whole paragraph\n$anchorStart = $rng.Start\n\n# Resolve the absolute paragraph index that starts at $anchorStart.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $anchorStart) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve the anchor paragraph's index\"\n}\n\n# Delete the 4 paragraphs right after the anchor as a single range.\n$deleteCount = 4\n$delStart = $d.Paragraphs.Item($anchorIndex + 1).Range.Start\n$delEnd = $d.Paragraphs.Item($anchorIndex + $deleteCount).Range.End\n$d.Range($delStart, $delEnd).Delete()\n"}
